# Updated assay metadata for imaging files.
#
# The "formatted.tar.gz" archive entry (row 2 of the "files" sheet) is no
# longer part of the file manifest, so it is removed. Deleting the entire
# row shifts every subsequent file row up by one (the shared-string table
# is pruned of the now-unused "formatted.tar.gz" / ".gz" strings as a
# natural consequence). This mirrors the interactive workflow of selecting
# the row on the "files" tab and deleting it, which is why "files" ends up
# as the active sheet/tab afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("files")

$ws.Activate() | Out-Null
$ws.Rows.Item(2).Select() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null
